$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.591.79"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.84%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.469.39"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -6.55%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "553.47"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.68%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.26"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.16%  "

$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.596"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -5.24%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.465.00"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -6.69%  "

$ws.Range("E10").Value = "  -9.04%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.47"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -6.13%  "

$ws.Range("E12").Value = "  -1.62%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.356"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -7.48%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.35"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -8.38%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.918.05"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -6.48%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000168"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -9.67%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.508.72"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.82%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.471.44"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -6.62%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.17"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -8.25%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.21"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -7.51%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.22"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.70%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "321.24"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -7.04%  "

$ws.Range("E23").Value = "  -0.06%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.91"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.26%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "64.38"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.73%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0₃0998"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -11.67%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "564.08"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.41%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.609.09"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.86%  "

$ws.Range("E29").Value = "  -8.12%  "

$ws.Range("E30").Value = "  +0.15%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.34"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -10.62%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.79"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.04%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.150"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -7.18%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.92"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -7.00%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.61"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -7.67%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.95"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -10.36%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.90"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -10.70%  "

$ws.Range("E38").Value = "  +0.03%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.381"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.27%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.51"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.45%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "144.40"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.05%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.76"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -7.99%  "

$ws.Range("E43").Value = "  +0.05%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.45"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.48%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "40.60"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.82%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "147.71"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -9.40%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.63"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -7.26%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "21.73"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -10.30%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0540"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -8.48%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.596"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.25%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0940"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -6.28%  "
